$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs","Tgfb2","Tgfbr3","ECs",3,1,2.441378666666667,7.324135999999999,0.1119936059016048,0.1119936059016048,3,1,32.56613966666666,97.698419,0.4203775945150413,0.4203775945150412,79.5062786378871,715.556507740984,0.04707960264998216,0.04707960264998214),
  @("ECs","Tgfb2","Tgfbr3","FAPs",3,1,2.441378666666667,7.324135999999999,0.1119936059016048,0.1119936059016048,3,1,40.85529333333333,122.56588,0.5273775188114274,0.5273775188114271,99.74324156440888,897.6891740796799,0.05906291000313317,0.05906291000313314),
  @("ECs","Tgfb2","Tgfbr3","M1",3,1,2.441378666666667,7.324135999999999,0.1119936059016048,0.1119936059016048,2,0.6666666666666666,0.03779266666666666,0.113378,0.0004878438300104565,0.0004878438300104564,0.09226621015644443,0.8303958914079999,0.00005463538963972055,0.00005463538963972053),
  @("ECs","Tgfb2","Tgfbr3","M2",3,1,2.441378666666667,7.324135999999999,0.1119936059016048,0.1119936059016048,3,1,0.07240166666666666,0.217205,0.0009345915353721286,0.0009345915353721284,0.1767598844311111,1.59083895988,0.0001046682760914419,0.0001046682760914419),
  @("ECs","Tgfb2","Tgfbr3","sCs",3,1,2.441378666666667,7.324135999999999,0.1119936059016048,0.1119936059016048,3,1,3.937153333333333,11.81146,0.05082245130814891,0.05082245130814889,9.612082155395555,86.50873939856,0.005691789582758328,0.005691789582758326),
  @("FAPs","Tgfb2","Tgfbr3","ECs",3,1,12.88577866666667,38.657336,0.5911106037886134,0.5911106037886134,3,1,32.56613966666666,97.698419,0.4203775945150413,0.4203775945150412,419.6400677724204,3776.760609951784,0.2484896537129909,0.2484896537129909),
  @("FAPs","Tgfb2","Tgfbr3","FAPs",3,1,12.88577866666667,38.657336,0.5911106037886134,0.5911106037886134,3,1,40.85529333333333,122.56588,0.5273775188114274,0.5273775188114271,526.4522672550755,4738.07040529568,0.3117384435691636,0.3117384435691635),
  @("FAPs","Tgfb2","Tgfbr3","M1",3,1,12.88577866666667,38.657336,0.5911106037886134,0.5911106037886134,2,0.6666666666666666,0.03779266666666666,0.113378,0.0004878438300104565,0.0004878438300104564,0.4869879378897777,4.382891441008,0.0002883696609120306,0.0002883696609120305),
  @("FAPs","Tgfb2","Tgfbr3","M2",3,1,12.88577866666667,38.657336,0.5911106037886134,0.5911106037886134,3,1,0.07240166666666666,0.217205,0.0009345915353721286,0.0009345915353721284,0.9329518517644443,8.39656666588,0.0005524469667695461,0.000552446966769546),
  @("FAPs","Tgfb2","Tgfbr3","sCs",3,1,12.88577866666667,38.657336,0.5911106037886134,0.5911106037886134,3,1,3.937153333333333,11.81146,0.05082245130814891,0.05082245130814889,50.73328643006222,456.59957787056,0.0300416898787773,0.03004168987877729),
  @("sCs","Tgfb2","Tgfbr3","ECs",3,1,6.472111000000001,19.416333,0.2968957903097819,0.2968957903097818,3,1,32.56613966666666,97.698419,0.4203775945150413,0.4203775945150412,210.7716707641697,1896.945036877527,0.1248083381520682,0.1248083381520682),
  @("sCs","Tgfb2","Tgfbr3","FAPs",3,1,6.472111000000001,19.416333,0.2968957903097819,0.2968957903097818,3,1,40.85529333333333,122.56588,0.5273775188114274,0.5273775188114271,264.4199933908934,2379.77994051804,0.1565761652391306,0.1565761652391305),
  @("sCs","Tgfb2","Tgfbr3","M1",3,1,6.472111000000001,19.416333,0.2968957903097819,0.2968957903097818,2,0.6666666666666666,0.03779266666666666,0.113378,0.0004878438300104565,0.0004878438300104564,0.2445983336526667,2.201385002874,0.0001448387794587053,0.0001448387794587053),
  @("sCs","Tgfb2","Tgfbr3","M2",3,1,6.472111000000001,19.416333,0.2968957903097819,0.2968957903097818,3,1,0.07240166666666666,0.217205,0.0009345915353721286,0.0009345915353721284,0.4685916232516666,4.217324609265,0.0002774762925111406,0.0002774762925111404),
  @("sCs","Tgfb2","Tgfbr3","sCs",3,1,6.472111000000001,19.416333,0.2968957903097819,0.2968957903097818,3,1,3.937153333333333,11.81146,0.05082245130814891,0.05082245130814889,25.48169339735334,229.33524057618,0.01508897184661328,0.01508897184661327)
)

$r = 2
foreach ($row in $rows) {
  for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c-1]
  }
  $r = $r + 1
}

Write-Output "done"